$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Bills"

$chartObj = $ws1.Shapes.AddChart2(227, 51)
$chart = $chartObj.Chart
$chart.ChartTitle.Text = "Hello world"
$chart.SetSourceData($ws1.Range("A1:B5"))
Write-Host "ok"
